# Updated cryptos list on Tue Sep 24 22:48:40 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row, and swaps the dogwifhat / VeChain rows (48 <-> 49) to match
# the latest coinranking.com snapshot.
#
# Note: several new Price values look like plain numbers (e.g. "608.83",
# "1.00") but must stay text (matching the original "Price" column, which
# is formatted as text, e.g. "64.364.29"). A leading "'" forces Excel to
# store them as text instead of auto-coercing to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.300.68"
$ws.Range("E2").Value = "  +1.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.655.14"
$ws.Range("E3").Value = "  +0.33%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'608.83"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'152.28"
$ws.Range("E6").Value = "  +5.72%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.83%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.70%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "'0.389"
$ws.Range("E10").Value = "  +6.59%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.62"
$ws.Range("E11").Value = "  -0.23%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.71%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'27.96"
$ws.Range("E13").Value = "  +2.30%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.131.06"
$ws.Range("E14").Value = "  +0.29%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "64.168.78"
$ws.Range("E15").Value = "  +1.68%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +2.19%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.681.26"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +6.30%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +4.59%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'346.90"
$ws.Range("E20").Value = "  +1.05%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.23%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.11%  "

# Row 23 - LEO
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "  +0.19%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'66.75"
$ws.Range("E24").Value = "  -0.58%  "

# Row 25 - SuiNetwork
$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = "  +15.17%  "

# Row 26 - Fetch.AI
$ws.Range("D26").Value = "'1.73"
$ws.Range("E26").Value = "  +4.90%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'9.42"
$ws.Range("E27").Value = "  +8.80%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "'8.25"
$ws.Range("E28").Value = "  +4.51%  "

# Row 29 - Bittensor
$ws.Range("D29").Value = "'555.63"
$ws.Range("E29").Value = "  +1.56%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +0.18%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.08%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +1.54%  "

# Row 33 - PEPE
$ws.Range("D33").Value = "0.0₃0861"
$ws.Range("E33").Value = "  +6.27%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  -1.02%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  +4.89%  "

# Row 36 - Monero
$ws.Range("D36").Value = "'168.45"
$ws.Range("E36").Value = "  -2.26%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("D37").Value = "'0.408"
$ws.Range("E37").Value = "  +0.60%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.01%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +6.63%  "

# Row 40 - EthereumClassic
$ws.Range("D40").Value = "'19.37"
$ws.Range("E40").Value = "  +1.27%  "

# Row 41 - USDe
$ws.Range("E41").Value = "  +0.10%  "

# Row 42 - Aave
$ws.Range("D42").Value = "'166.82"
$ws.Range("E42").Value = "  -3.04%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'40.36"
$ws.Range("E43").Value = "  +0.68%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +2.82%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "'0.0577"
$ws.Range("E45").Value = "  -0.18%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'22.03"
$ws.Range("E46").Value = "  -1.44%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "'0.630"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48 - was VeChain, now dogwifhat
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = "  +15.12%  "

# Row 49 - was dogwifhat, now VeChain
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0246"
$ws.Range("E49").Value = "  +2.66%  "

# Row 50 - Stellar
$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  +0.32%  "

# Row 51 - EnergySwap
$ws.Range("E51").Value = "  +1.53%  "
